$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.042870879173279
$ws.Range("B1").Value = 1.302420139312744
$ws.Range("D1").Value = 1.683874130249023
$ws.Range("E1").Value = 1.014582276344299
